$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 800
$ws.Range("B3").Value = 400
$ws.Range("B4").Value = 400
$ws.Range("B5").Value = 100
$ws.Range("B6").Value = 50
$ws.Range("B7").Value = 500
$ws.Range("B8").Value = 228
